$wb = $excel.ActiveWorkbook

# New identifiers for this handoff/handback cycle
$oldGuid = "bbccf146-a698-47e6-80fd-e51a15bf10b4"
$newGuid = "47f13380-eab1-49b2-a11f-c55e9360459a"
$oldHash = "6b4bf3f50d20549ae0b8352b9ae11dac1f34cdc5"
$newHash = "ed27e36f28d8bceca47b46c2b7cf0a71474f6f9b"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-24 16:59:54"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-24 16:59:50"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"

# --- Update hyperlink display text (targets unchanged) ---
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
